$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "이노그리드" IPO entry moves up in the schedule: insert a fresh row
# for it (with updated demand-forecast date / price band / raised amount /
# lead underwriter) right after "이엔셀" (row 3), and remove the old entry
# further down the table.

# 1) Find the row that currently holds the old "이노그리드" entry.
$oldRow = $null
for ($r = 2; $r -le $ws.UsedRange.Rows.Count; $r++) {
    if ($ws.Cells.Item($r, 1).Value2 -eq "이노그리드") {
        $oldRow = $r
        break
    }
}

# 2) Insert a new blank row right before row 4 (after "이엔셀"), pushing
#    everything from row 4 down (including the old "이노그리드" row).
$ws.Rows.Item(4).Insert(-4121) | Out-Null   # xlShiftDown = -4121

# 3) Populate the new row 4 with the updated "이노그리드" data.
$newRow = 4
$ws.Cells.Item($newRow, 1).Value2 = "이노그리드"
$ws.Cells.Item($newRow, 2).Value2 = "2024.06.13~06.19"
$ws.Cells.Item($newRow, 3).Value2 = "29,000~35,000"
$ws.Cells.Item($newRow, 4).Value2 = "-"
$ws.Cells.Item($newRow, 5).Value2 = 17400
$ws.Cells.Item($newRow, 6).Value2 = "한국투자증권"

# 4) Delete the old "이노그리드" row (it shifted down by one after the insert).
$ws.Rows.Item($oldRow + 1).Delete(-4162) | Out-Null   # xlShiftUp = -4162
